$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.552.50'
$ws.Range('E2').Value = '  +6.71%  '
$ws.Range('D3').Value = '2.582.92'
$ws.Range('E3').Value = '  +8.70%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '504.74'
$ws.Range('E5').Value = '  +5.42%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '157.18'
$ws.Range('E6').Value = '  +6.70%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.623'
$ws.Range('E7').Value = '  +23.90%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').Value = '2.579.31'
$ws.Range('E9').Value = '  +8.53%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.17'
$ws.Range('E10').Value = '  +13.83%  '
$ws.Range('E11').Value = '  +5.68%  '
$ws.Range('E12').Value = '  +5.41%  '
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').Value = '3.021.58'
$ws.Range('E14').Value = '  +8.40%  '
$ws.Range('D15').Value = '59.486.61'
$ws.Range('E15').Value = '  +6.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '21.89'
$ws.Range('E16').Value = '  +7.65%  '
$ws.Range('E17').Value = '  +3.93%  '
$ws.Range('D18').Value = '2.580.58'
$ws.Range('E18').Value = '  +8.89%  '
$ws.Range('E19').Value = '  +2.90%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '334.86'
$ws.Range('E20').Value = '  +6.11%  '
$ws.Range('E21').Value = '  +7.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.06'
$ws.Range('E22').Value = '  +7.04%  '
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '60.40'
$ws.Range('E24').Value = '  +6.44%  '
$ws.Range('E25').Value = '  +4.74%  '
$ws.Range('E26').Value = '  +7.48%  '
$ws.Range('D27').Value = '2.676.62'
$ws.Range('E27').Value = '  +7.85%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.17%  '
$ws.Range('E29').Value = '  +2.71%  '
$ws.Range('D30').Value = '0.0₃0822'
$ws.Range('E30').Value = '  +6.44%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '156.67'
$ws.Range('E32').Value = '  +6.12%  '
$ws.Range('E33').Value = '  +7.74%  '
$ws.Range('E34').Value = '  +5.30%  '
$ws.Range('E35').Value = '  +8.10%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.92'
$ws.Range('E36').Value = '  +9.64%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.19'
$ws.Range('E37').Value = '  +7.93%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.860'
$ws.Range('E38').Value = '  +2.89%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.79'
$ws.Range('E39').Value = '  +12.15%  '
$ws.Range('E40').Value = '  +7.43%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '34.95'
$ws.Range('E41').Value = '  +4.59%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '291.76'
$ws.Range('E42').Value = '  +14.25%  '
$ws.Range('E43').Value = '  +7.40%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.625'
$ws.Range('E44').Value = '  +7.37%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0561'
$ws.Range('E45').Value = '  +4.42%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '19.34'
$ws.Range('E47').Value = '  +14.65%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0237'
$ws.Range('E48').Value = '  +6.32%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.76'
$ws.Range('E49').Value = '  +5.09%  '
$ws.Range('E50').Value = '  +1.04%  '
$ws.Range('E51').Value = '  +12.35%  '
